$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The matrix grew from 3x3 to a 4x4 "gauss iteration" matrix and the
# cell number format changed from an integer format (#,##0) to a
# 2-decimal format (#,##0.00). Grab the new format onto A1 first, then
# fan it out (via copy / paste-special-formats) to the whole new A1:D4
# block so every cell ends up sharing one consistent style - mirroring
# how the original A1:C3 block shared a single style.
$ws.Range("A1").NumberFormat = "#,##0.00"
$ws.Range("A1").Copy()
$ws.Range("A1:D4").PasteSpecial(-4122)

$values = @(
    @(10, 5, -2, 4),
    @(3, 9, -1, 2),
    @(0, 4, 15, 5),
    @(-1, 2, 3, 7)
)

for ($r = 1; $r -le 4; $r++) {
    for ($c = 1; $c -le 4; $c++) {
        $ws.Cells.Item($r, $c).Value = $values[$r - 1][$c - 1]
    }
    $ws.Rows.Item($r).RowHeight = 19.5
}

# New column D mirrors the width of the existing bestFit columns A:C.
$ws.Columns("D").ColumnWidth = 12.67
